$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: find the (first remaining) occurrence of $searchText in the whole
# document, replace it with the concatenation of $parts, then force the
# result to be split into one run per element of $parts by toggling Bold
# on/off across the tail of the newly-inserted text at each split boundary
# (a no-op formatting change that nonetheless forces Word's run-coalescing
# engine to keep the pieces as separate <w:r> elements, matching the way
# real Word leaves behind multiple runs after a sequence of small edits).
# ---------------------------------------------------------------------------
function ReplaceAndSplit($searchText, $parts) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND: $searchText"
        return $false
    }
    $start = $rng.Start
    $end = $rng.End
    $newText = [string]::Join("", $parts)
    $full = $d.Range($start, $end)
    $full.Text = $newText
    $newEnd = $start + $newText.Length
    $offset = 0
    for ($i = 0; $i -lt ($parts.Length - 1); $i++) {
        $offset = $offset + $parts[$i].Length
        $pos = $start + $offset
        if ($pos -lt $newEnd) {
            $r = $d.Range($pos, $newEnd)
            $r.Bold = 1
            $r.Bold = 0
        }
    }
    return $true
}

# Apply ReplaceAndSplit exactly $count times (each table text below shows up
# twice in this document — once per ranova table).
function ReplaceAndSplitN($searchText, $parts, $count) {
    for ($k = 0; $k -lt $count; $k++) {
        $ok = ReplaceAndSplit $searchText $parts
        if (-not $ok) {
            Write-Host "STOPPED EARLY: $searchText (iteration $k)"
            break
        }
    }
}

# ---------------------------------------------------------------------------
# 1) Intro paragraph: "Values for last year included (e.g. ... below)"
# ---------------------------------------------------------------------------
ReplaceAndSplitN "Values for last year included (e.g. relative growth rate measured 2019-2021; 2021 values are below)" @(
    "Values for last year included (",
    "e.g.",
    " relative growth rate measured 2019-2021; 2021 values are below)"
) 1

# ---------------------------------------------------------------------------
# 2) Table row labels (each appears twice — once per table)
# ---------------------------------------------------------------------------
ReplaceAndSplitN "Height, before flowering" @("Height", " b", "efore flowering") 2
ReplaceAndSplitN "Height, after flowering" @("Height", " a", "fter flowering") 2

ReplaceAndSplitN "Ramets, before flowering*" @("Ramets", " b", "efore flowering*") 2
ReplaceAndSplitN "Ramets, after flowering*" @("Ramets", " a", "fter flowering*") 2

ReplaceAndSplitN "Danaus plexippus" @("D.", " plexippus") 2
ReplaceAndSplitN "Labidomera clivicollis " @("L.", " clivicollis ") 2
ReplaceAndSplitN "Liriomyza asclepiadis" @("L.", " asclepiadis") 2

ReplaceAndSplitN "Herbivory, before flowering: Binary*" @("Herbivory", " b", "efore flowering", " (binary)", "*") 2
ReplaceAndSplitN "Herbivory, before flowering: Quantitative" @("Herbivory", " b", "efore flowering", " (quantitative)") 2
ReplaceAndSplitN "Herbivory, after flowering: Binary*" @("Herbivory", " a", "fter flowering", " (binary)", "*") 2
ReplaceAndSplitN "Herbivory, after flowering: Quantitative" @("Herbivory", " a", "fter flowering", " (quantitative)") 2

ReplaceAndSplitN "Weevil damage: Binary*" @("Weevil damage", " (binary)", "*") 2
ReplaceAndSplitN "Weevil damage: Quantitative" @("Weevil damage", " (quantitative)") 2

ReplaceAndSplitN "Mean flowers per inflorescence*" @("Fl", "owers ", "per in", "florescence*") 2

# ---------------------------------------------------------------------------
# 3) Footnote paragraph: PVE calculation explanation (get_variance / VarCorr)
# ---------------------------------------------------------------------------
ReplaceAndSplitN "*Variables were analyzed with generalized linear mixed models. PVE was calculated as: random effect variance/(random effect variance + residual variance) with the" @(
    "*Variables were analyzed with generalized linear mixed models. PVE was calculated as: random effect variance",
    "/(",
    "random effect variance + residual variance) with the"
) 1

ReplaceAndSplitN "get_variance()" @("get_variance", "()") 1

ReplaceAndSplitN "random effect variance/(random effect variance + residual variance) with the " @(
    "random effect variance",
    "/(",
    "random effect variance + residual variance) with the "
) 1

ReplaceAndSplitN "VarCorr()" @("VarCorr", "()") 1
